# Refresh the cryptos price table (columns B-E, rows 2-51) with the latest
# scraped values from coinranking.com, as produced by the scheduled
# GitHub Actions job.
#
# Note: several "Price" values (column D) look like plain numbers to Excel
# (e.g. "1.001", "317.38"), which would otherwise be auto-converted to a
# numeric type and lose their exact textual formatting/trailing zeros.
# Prefixing them with a leading single-quote forces Excel to keep them as
# literal text, matching the original inline-string cell content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''28.207.38'
$ws.Range('E2').Value = '  -2.36%  '

$ws.Range('D3').Value = '''1.805.73'
$ws.Range('E3').Value = '  -0.33%  '

$ws.Range('E4').Value = '  +0.26%  '

$ws.Range('D5').Value = '''317.38'
$ws.Range('E5').Value = '  +0.66%  '

$ws.Range('D6').Value = '''1.001'
$ws.Range('E6').Value = '  +0.48%  '

$ws.Range('D7').Value = '''0.5308'
$ws.Range('E7').Value = '  -2.18%  '

$ws.Range('E8').Value = '  -1.27%  '

$ws.Range('D9').Value = '''0.07492'
$ws.Range('E9').Value = '  -1.39%  '

$ws.Range('D10').Value = '''41.98'
$ws.Range('E10').Value = '  -1.83%  '

$ws.Range('D11').Value = '''1.100'
$ws.Range('E11').Value = '  -2.61%  '

$ws.Range('D12').Value = '''1.001'
$ws.Range('E12').Value = '  +0.28%  '

$ws.Range('D13').Value = '''6.223'
$ws.Range('E13').Value = '  -0.21%  '

$ws.Range('E14').Value = '  -3.59%  '

$ws.Range('D15').Value = '''7.362'
$ws.Range('E15').Value = '  -1.63%  '

$ws.Range('D16').Value = '''1.803.09'
$ws.Range('E16').Value = '  -0.08%  '

$ws.Range('D17').Value = '''89.72'
$ws.Range('E17').Value = '  -2.73%  '

$ws.Range('D18').Value = '''0.00001070'
$ws.Range('E18').Value = '  -0.20%  '

$ws.Range('E19').Value = '  +0.93%  '

$ws.Range('D20').Value = '''1.0000'
$ws.Range('E20').Value = '  +0.47%  '

$ws.Range('D21').Value = '''17.26'
$ws.Range('E21').Value = '  -0.62%  '

$ws.Range('D22').Value = '''5.925'
$ws.Range('E22').Value = '  -1.27%  '

$ws.Range('D23').Value = '''28.231.13'
$ws.Range('E23').Value = '  -2.34%  '

$ws.Range('D24').Value = '''11.19'
$ws.Range('E24').Value = '  -2.44%  '

$ws.Range('D25').Value = '''2.090'
$ws.Range('E25').Value = '  -1.81%  '

$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '''20.50'
$ws.Range('E26').Value = '  -1.14%  '

$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = '''155.55'
$ws.Range('E27').Value = '  -4.23%  '

$ws.Range('D28').Value = '''2.011.88'
$ws.Range('E28').Value = '  -0.21%  '

$ws.Range('E29').Value = '  -3.94%  '

$ws.Range('D30').Value = '''122.12'
$ws.Range('E30').Value = '  -1.94%  '

$ws.Range('E31').Value = '  -2.38%  '

$ws.Range('D32').Value = '''0.1099'
$ws.Range('E32').Value = '  +7.67%  '

$ws.Range('D33').Value = '''5.593'
$ws.Range('E33').Value = '  -3.80%  '

$ws.Range('D34').Value = '''3.627'
$ws.Range('E34').Value = '  -1.23%  '

$ws.Range('D35').Value = '''0.07234'
$ws.Range('E35').Value = '  +8.92%  '

$ws.Range('E36').Value = '  -4.08%  '

$ws.Range('E37').Value = '  -1.57%  '

$ws.Range('D38').Value = '''5.100'
$ws.Range('E38').Value = '  -1.11%  '

$ws.Range('D39').Value = '''8.518'
$ws.Range('E39').Value = '  -1.62%  '

$ws.Range('D40').Value = '''0.6177'
$ws.Range('E40').Value = '  -3.46%  '

$ws.Range('D41').Value = '''11.14'
$ws.Range('E41').Value = '  -4.70%  '

$ws.Range('E42').Value = '  -4.03%  '

$ws.Range('E43').Value = '  +1.94%  '

$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '''13.34'
$ws.Range('E44').Value = '  -2.91%  '

$ws.Range('B45').Value = 'PancakeSwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D45').Value = '''3.681'
$ws.Range('E45').Value = '  -0.11%  '

$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '''0.5775'
$ws.Range('E46').Value = '  -3.89%  '

$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').Value = '''125.90'
$ws.Range('E47').Value = '  -0.14%  '

$ws.Range('B48').Value = 'EOS'
$ws.Range('C48').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D48').Value = '''1.195'
$ws.Range('E48').Value = '  +2.11%  '

$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '''1.928'
$ws.Range('E49').Value = '  -4.18%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '''0.06824'
$ws.Range('E50').Value = '  -2.33%  '

$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '''72.19'
$ws.Range('E51').Value = '  -1.84%  '
